# Commit atualizacao tesouro dia 10/06
# Fills in the next two trading-day columns (07/jun and 10/jun) of the
# "GRAFICO" sheet: purchase-title rate/price (row 4) and sale-title
# rate/price (row 5), plus their date header (row 3, columns L:O).
# Downstream formulas (rows 6-7) and the charts that read this range
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GRAFICO")

# Row 3: date headers for 07/jun (L,M) and 10/jun (N,O)
$ws.Range("L3").Value = 43623
$ws.Range("M3").Value = 43623
$ws.Range("N3").Value = 43626
$ws.Range("O3").Value = 43626

# Row 5 (TITULO VENDA line) picks up the border formatting already used
# by the previous day's rate cells (H5/J5) when the new values are typed
# in, so copy that formatting across before setting the values.
$ws.Range("H5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("J5").Copy()
$ws.Range("N5").PasteSpecial(-4122)

# Row 4 (TITULO COMPRA line): rate / price for 07/jun and 10/jun
$ws.Range("L4").Value = 4.05
$ws.Range("M4").Value = 1718.57
$ws.Range("N4").Value = 4.05
$ws.Range("O4").Value = 1718.91

# Row 5 (TITULO VENDA line): rate / price for 07/jun and 10/jun
$ws.Range("L5").Value = 3.93
$ws.Range("M5").Value = 1750.35
$ws.Range("N5").Value = 3.93
$ws.Range("O5").Value = 1750.69
